$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 338, shifting existing rows 338-368 down to 339-369
$ws.Rows.Item(338).Insert()

# Copy the date number format from the row above (so the new date cell
# matches the existing "$/paquete ..." rows' date formatting)
$ws.Range("D337").Copy()
$ws.Range("D338").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Populate the new row with the weekly price update
$ws.Range("A338").Value = 8
$ws.Range("B338").Value = "Terminal La Palmera de La Serena"
$ws.Range("C338").Value = "Coquimbo"
$ws.Range("D338").Value = Get-Date -Year 2023 -Month 9 -Day 25 -Hour 0 -Minute 0 -Second 0
$ws.Range("E338").Value = 4
$ws.Range("F338").Value = 100112037
$ws.Range("G338").Value = "Cebollín"
$ws.Range("H338").Value = "Sin especificar"
$ws.Range("I338").Value = "Primera"
$ws.Range("J338").Value = 1800
$ws.Range("K338").Value = 1000
$ws.Range("L338").Value = 1200
$ws.Range("M338").Value = 1100
$ws.Range("N338").Value = "$/paquete 6 unidades"
$ws.Range("O338").Value = "Provincia del Elquí"
$ws.Range("P338").Value = 183
$ws.Range("Q338").Value = 6
$ws.Range("R338").Value = "Hortaliza"
